# Auto-generated edit script replicating the WORDS.xlsx diff:
# - adds a new block of 30 GRE words (set 11) to Sheet2, rows 332-361
# - re-labels the previously-merged second half of set 9 (rows 302-331) as set 10
# - grows the autofilter range and the Sheet2!_FilterDatabase defined name to A1:E361

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# --- Rows 302-331 were labelled "9" (continuation of the 272-331 block);
# --- the edit splits that block so 302-331 becomes its own "10" group.
for ($r = 302; $r -le 331; $r++) {
  $ws.Cells.Item($r, 1).Value = 10
}

# --- Append the new "11" word group: columns A (set no.), B (serial no.),
# --- C (word), D (meaning), E (example sentence).
$r = 332
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 331
$ws.Cells.Item($r, 3).Value = "Acrimonious"
$ws.Cells.Item($r, 4).Value = "argumentative, threatening, or aggressive"
$ws.Cells.Item($r, 5).Value = "The debate between the two candidates became acrimonious, with personal attacks overshadowing the issues."
$ws.Cells.Item($r, 3).VerticalAlignment = -4108
$ws.Cells.Item($r, 5).VerticalAlignment = -4108

$r = 333
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 332
$ws.Cells.Item($r, 3).Value = "Belligerent"
$ws.Cells.Item($r, 4).Value = "argumentative, threatening, or aggressive"
$ws.Cells.Item($r, 5).Value = "His belligerent tone during the negotiation only made it harder to reach a peaceful resolution."
$ws.Cells.Item($r, 3).VerticalAlignment = -4108
$ws.Cells.Item($r, 5).VerticalAlignment = -4108

$r = 334
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 333
$ws.Cells.Item($r, 3).Value = "Beneficent"
$ws.Cells.Item($r, 4).Value = "kind and charitable, producing good results or benefits"
$ws.Cells.Item($r, 5).Value = "The beneficent philanthropist donated millions to help underprivileged children access education."
$ws.Cells.Item($r, 3).VerticalAlignment = -4108
$ws.Cells.Item($r, 5).VerticalAlignment = -4108

$r = 335
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 334
$ws.Cells.Item($r, 3).Value = "Canny"
$ws.Cells.Item($r, 4).Value = "shrewd and cautious"
$ws.Cells.Item($r, 5).Value = "Her canny investments in the stock market earned her significant returns despite the economic downturn."
$ws.Cells.Item($r, 3).VerticalAlignment = -4108
$ws.Cells.Item($r, 5).VerticalAlignment = -4108

$r = 336
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 335
$ws.Cells.Item($r, 3).Value = "Cavalier"
$ws.Cells.Item($r, 4).Value = "exhibiting an overly casual attitude; unconcerned"
$ws.Cells.Item($r, 5).Value = "His cavalier attitude toward deadlines frustrated his team, who were working hard to meet the schedule."
$ws.Cells.Item($r, 3).VerticalAlignment = -4108
$ws.Cells.Item($r, 5).VerticalAlignment = -4108

$r = 337
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 336
$ws.Cells.Item($r, 3).Value = "Distressed"
$ws.Cells.Item($r, 4).Value = "experiencing nervousness, irritation, or sadness"
$ws.Cells.Item($r, 5).Value = "She became distressed after learning about the unexpected delay in her project’s approval."
$ws.Cells.Item($r, 3).VerticalAlignment = -4108
$ws.Cells.Item($r, 5).VerticalAlignment = -4108

$r = 338
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 337
$ws.Cells.Item($r, 3).Value = "Dwindling"
$ws.Cells.Item($r, 4).Value = "decreasing steadily in size or quantity"
$ws.Cells.Item($r, 5).Value = "The company’s dwindling profits forced them to reconsider their business strategy."
$ws.Cells.Item($r, 3).VerticalAlignment = -4108
$ws.Cells.Item($r, 5).VerticalAlignment = -4108

$r = 339
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 338
$ws.Cells.Item($r, 3).Value = "Eclipse"
$ws.Cells.Item($r, 4).Value = "overshadow or surpass"
$ws.Cells.Item($r, 5).Value = "His groundbreaking research managed to eclipse all previous studies on the subject."
$ws.Cells.Item($r, 3).VerticalAlignment = -4108
$ws.Cells.Item($r, 5).VerticalAlignment = -4108

$r = 340
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 339
$ws.Cells.Item($r, 3).Value = "Encyclopedic"
$ws.Cells.Item($r, 4).Value = "thorough; comprehensive in scope"
$ws.Cells.Item($r, 5).Value = "Her encyclopedic knowledge of history made her an invaluable resource for the research team."
$ws.Cells.Item($r, 3).VerticalAlignment = -4108
$ws.Cells.Item($r, 5).VerticalAlignment = -4108

$r = 341
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 340
$ws.Cells.Item($r, 3).Value = "Exacerbate"
$ws.Cells.Item($r, 4).Value = "make a situation or condition worse"
$ws.Cells.Item($r, 5).Value = "The manager’s refusal to address the complaints only served to exacerbate the tensions among employees."
$ws.Cells.Item($r, 3).VerticalAlignment = -4108
$ws.Cells.Item($r, 5).VerticalAlignment = -4108

$r = 342
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 341
$ws.Cells.Item($r, 3).Value = "Exasperated"
$ws.Cells.Item($r, 4).Value = "severely irritated or angry"
$ws.Cells.Item($r, 5).Value = "She was exasperated by his repeated excuses for missing important meetings."
$ws.Cells.Item($r, 3).VerticalAlignment = -4108
$ws.Cells.Item($r, 5).VerticalAlignment = -4108

$r = 343
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 342
$ws.Cells.Item($r, 3).Value = "Fungible"
$ws.Cells.Item($r, 4).Value = "interchangeable for another item of a similar type"
$ws.Cells.Item($r, 5).Value = "In the commodities market, gold and silver are considered fungible assets, easily traded for one another."
$ws.Cells.Item($r, 3).VerticalAlignment = -4108
$ws.Cells.Item($r, 5).VerticalAlignment = -4108

$r = 344
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 343
$ws.Cells.Item($r, 3).Value = "Hackneyed"
$ws.Cells.Item($r, 4).Value = "lacking creativity; derivative and overused"
$ws.Cells.Item($r, 5).Value = "The movie’s hackneyed plot failed to engage the audience, as it felt predictable and overused."
$ws.Cells.Item($r, 3).VerticalAlignment = -4108
$ws.Cells.Item($r, 5).VerticalAlignment = -4108

$r = 345
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 344
$ws.Cells.Item($r, 3).Value = "Incongruous"
$ws.Cells.Item($r, 4).Value = "out of place; lacking harmony"
$ws.Cells.Item($r, 5).Value = "The modern design of the building looked incongruous in the historic neighborhood."
$ws.Cells.Item($r, 3).VerticalAlignment = -4108
$ws.Cells.Item($r, 5).VerticalAlignment = -4108

$r = 346
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 345
$ws.Cells.Item($r, 3).Value = "Interchangeable"
$ws.Cells.Item($r, 4).Value = "capable of being used in place of each other"
$ws.Cells.Item($r, 5).Value = "The two terms are not interchangeable, as each has a distinct meaning in legal contexts."
$ws.Cells.Item($r, 3).VerticalAlignment = -4108
$ws.Cells.Item($r, 5).VerticalAlignment = -4108

$r = 347
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 346
$ws.Cells.Item($r, 3).Value = "Laconic"
$ws.Cells.Item($r, 4).Value = "economical with words; quiet and reserved"
$ws.Cells.Item($r, 5).Value = "His laconic reply, consisting of just a single word, left everyone in the room confused."
$ws.Cells.Item($r, 3).VerticalAlignment = -4108
$ws.Cells.Item($r, 5).VerticalAlignment = -4108

$r = 348
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 347
$ws.Cells.Item($r, 3).Value = "Lucrative"
$ws.Cells.Item($r, 4).Value = "producing a significant amount of money"
$ws.Cells.Item($r, 5).Value = "The startup turned out to be a lucrative venture, generating substantial profits within its first year."
$ws.Cells.Item($r, 3).VerticalAlignment = -4108
$ws.Cells.Item($r, 5).VerticalAlignment = -4108

$r = 349
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 348
$ws.Cells.Item($r, 3).Value = "Magisterial"
$ws.Cells.Item($r, 4).Value = "possessing great authority"
$ws.Cells.Item($r, 5).Value = "The professor’s magisterial presence in the lecture hall commanded the attention of all the students."
$ws.Cells.Item($r, 3).VerticalAlignment = -4108
$ws.Cells.Item($r, 5).VerticalAlignment = -4108

$r = 350
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 349
$ws.Cells.Item($r, 3).Value = "Onerous"
$ws.Cells.Item($r, 4).Value = "involving a heavy burden; challenging or difficult"
$ws.Cells.Item($r, 5).Value = "The contract imposed onerous conditions that made it nearly impossible for the company to comply."
$ws.Cells.Item($r, 3).VerticalAlignment = -4108
$ws.Cells.Item($r, 5).VerticalAlignment = -4108

$r = 351
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 350
$ws.Cells.Item($r, 3).Value = "Opprobrium"
$ws.Cells.Item($r, 4).Value = "severe criticism or public shame"
$ws.Cells.Item($r, 5).Value = "The politician faced widespread opprobrium after his unethical actions were exposed to the public."
$ws.Cells.Item($r, 3).VerticalAlignment = -4108
$ws.Cells.Item($r, 5).VerticalAlignment = -4108

$r = 352
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 351
$ws.Cells.Item($r, 3).Value = "Parsimonious"
$ws.Cells.Item($r, 4).Value = "extremely frugal or stingy"
$ws.Cells.Item($r, 5).Value = "The parsimonious landlord refused to spend money on essential repairs for the apartment."
$ws.Cells.Item($r, 3).VerticalAlignment = -4108
$ws.Cells.Item($r, 5).VerticalAlignment = -4108

$r = 353
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 352
$ws.Cells.Item($r, 3).Value = "Peripheral"
$ws.Cells.Item($r, 4).Value = "of marginal importance; minor"
$ws.Cells.Item($r, 5).Value = "The discussion focused on peripheral issues rather than addressing the core problem at hand."
$ws.Cells.Item($r, 3).VerticalAlignment = -4108
$ws.Cells.Item($r, 5).VerticalAlignment = -4108

$r = 354
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 353
$ws.Cells.Item($r, 3).Value = "Provocative"
$ws.Cells.Item($r, 4).Value = "causing outrage or irritation, often intentionally"
$ws.Cells.Item($r, 5).Value = "His provocative remarks during the panel discussion sparked a heated debate among the participants."
$ws.Cells.Item($r, 3).VerticalAlignment = -4108
$ws.Cells.Item($r, 5).VerticalAlignment = -4108

$r = 355
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 354
$ws.Cells.Item($r, 3).Value = "Renounce"
$ws.Cells.Item($r, 4).Value = "reject or give up"
$ws.Cells.Item($r, 5).Value = "She decided to renounce her title as a sign of solidarity with the common people."
$ws.Cells.Item($r, 3).VerticalAlignment = -4108
$ws.Cells.Item($r, 5).VerticalAlignment = -4108

$r = 356
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 355
$ws.Cells.Item($r, 3).Value = "Tempestuous"
$ws.Cells.Item($r, 4).Value = "(of emotions or actions) unrestrained and turbulent"
$ws.Cells.Item($r, 5).Value = "Their tempestuous relationship was marked by frequent arguments and passionate reconciliations."
$ws.Cells.Item($r, 3).VerticalAlignment = -4108
$ws.Cells.Item($r, 5).VerticalAlignment = -4108

$r = 357
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 356
$ws.Cells.Item($r, 3).Value = "Tenable"
$ws.Cells.Item($r, 4).Value = "(of an idea or method) justifiable and rational"
$ws.Cells.Item($r, 5).Value = "The scientist presented a tenable argument, backed by extensive data, to support her hypothesis."
$ws.Cells.Item($r, 3).VerticalAlignment = -4108
$ws.Cells.Item($r, 5).VerticalAlignment = -4108

$r = 358
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 357
$ws.Cells.Item($r, 3).Value = "Transgression"
$ws.Cells.Item($r, 4).Value = "a violation of a law, rule, or social norm"
$ws.Cells.Item($r, 5).Value = "The student’s transgression of cheating on the exam resulted in severe disciplinary action."
$ws.Cells.Item($r, 3).VerticalAlignment = -4108
$ws.Cells.Item($r, 5).VerticalAlignment = -4108

$r = 359
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 358
$ws.Cells.Item($r, 3).Value = "Urbane"
$ws.Cells.Item($r, 4).Value = "refined, sophisticated, and courteous"
$ws.Cells.Item($r, 5).Value = "His urbane manners and polished speech made him a favorite at social gatherings."
$ws.Cells.Item($r, 3).VerticalAlignment = -4108
$ws.Cells.Item($r, 5).VerticalAlignment = -4108

$r = 360
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 359
$ws.Cells.Item($r, 3).Value = "Verisimilitude"
$ws.Cells.Item($r, 4).Value = "appearance or semblance of truth or reality; believability."
$ws.Cells.Item($r, 5).Value = "The novel’s vivid descriptions lent a sense of verisimilitude to the fictional world it portrayed."
$ws.Cells.Item($r, 3).VerticalAlignment = -4108
$ws.Cells.Item($r, 5).VerticalAlignment = -4108

$r = 361
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = 360
$ws.Cells.Item($r, 3).Value = "Vitiate"
$ws.Cells.Item($r, 4).Value = "impair or corrupt the quality of something"
$ws.Cells.Item($r, 5).Value = "The presence of bias in the study’s methodology could vitiate the reliability of its conclusions."
$ws.Cells.Item($r, 3).VerticalAlignment = -4108
$ws.Cells.Item($r, 5).VerticalAlignment = -4108

# --- Grow the autofilter to cover the new rows (toggle off/on since the
# --- sheet already has an active autofilter over A1:E301).
$ws.AutoFilterMode = $false
$ws.Range("A1:E361").AutoFilter()

# --- Grow the hidden _FilterDatabase defined name to match.
for ($i = 1; $i -le $wb.Names.Count; $i++) {
  $nm = $wb.Names.Item($i)
  if ($nm.Name -eq "Sheet2!_FilterDatabase") {
    $nm.RefersTo = "=Sheet2!`$A`$1:`$E`$361"
  }
}

# --- Leave the selection matching the saved view (the committed worksheet
# --- records activeCell="E351" / sqref="E351").
$ws.Range("E351").Select()

